$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-02-28"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 02-28)"

# Update February total (row 3) and grand Total row (row 14) for column I
$ws.Range("I3").Value = 141
$ws.Range("I14").Value = 300
